# Quarterly indexing esoteric bug-fix operation
# Fix the date serials in column A (rows 2-47): the dates were stored as the
# last instant of the quarter-end month (e.g. 40543.99999999999 ~ 2010-12-31)
# but should instead be the serial date value for the 15th of the month two
# months earlier (e.g. 40497 = 2010-11-15).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$dateValues = @{
    2  = 40497
    3  = 40678
    4  = 40862
    5  = 41044
    6  = 41228
    7  = 41409
    8  = 41593
    9  = 41774
    10 = 41958
    11 = 42139
    12 = 42323
    13 = 42505
    14 = 42689
    15 = 42870
    16 = 43054
    17 = 43146
    18 = 43235
    19 = 43327
    20 = 43419
    21 = 43511
    22 = 43600
    23 = 43692
    24 = 43784
    25 = 43876
    26 = 43966
    27 = 44058
    28 = 44150
    29 = 44242
    30 = 44331
    31 = 44423
    32 = 44515
    33 = 44607
    34 = 44696
    35 = 44788
    36 = 44880
    37 = 44972
    38 = 45061
    39 = 45153
    40 = 45245
    41 = 45337
    42 = 45427
    43 = 45519
    44 = 45611
    45 = 45703
    46 = 45792
    47 = 45884
}

foreach ($row in $dateValues.Keys) {
    $ws.Cells.Item($row, 1).Value = $dateValues[$row]
}
